$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Fifi" entry (row 5) entirely - shifts rows 6-8 up to 5-7
$ws.Rows.Item(5).Delete()

# Update the "path to pic" column (H) with renamed picture files
# (an underscore was inserted between the handler/kennel name and the dog name)
$ws.Range("H2").Value = "PicturesOrig\FolkeNoertemann_Arlo.jpg"
$ws.Range("H3").Value = "PicturesOrig\AnnieVanderlinck_Moss.jpg"
$ws.Range("H4").Value = "PicturesOrig\AnnieVanderlinck_Tweed.jpg"
$ws.Range("H5").Value = "PicturesOrig\FolkeNoertemannKinloch_Luke.jpg"
$ws.Range("H6").Value = "PicturesOrig\FolkeNoertemannKinloch_Heath.jpg"
$ws.Range("H7").Value = "PicturesOrig\FolkeNoertemann_Joe.jpg"

# Update the selected cell shown when the workbook was last saved
$ws.Range("H13").Select()
